{"js": "// Append new status-update paragraphs to the end of the document body,\n// matching the existing \"Helvetica Light\" / size 24 run formatting used\n// throughout the document.\n\nconst newLines = [\n  \"\",\n  \"Update: 2026-02-21 (Contacts)\",\n  \"- Added contact Segment mapping support (backend + frontend).\",\n  \"- Contacts now return segment name in list API for UI display.\",\n  \"- Added SegmentId to contact schema and upsert request.\",\n  \"- Outbound message send now auto-creates/updates contact by recipient number.\",\n  \"- Dashboard Contacts page now uses real segments API, supports create segment, and maps segment while creating contact.\",\n  \"- WABA Contacts page updated to map segment on add/edit and display segment column.\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfor (const line of newLines) {\n  const p = anchor.insertParagraph(\"\", \"After\");\n  p.insertText(line, \"Replace\");\n  p.font.name = \"Helvetica Light\";\n  p.font.size = 12;\n  anchor = p;\n}\n\nawait context.sync();\n", "ps1": "# Append new status-update paragraphs to the end of the document,\n# matching the existing \"Helvetica Light\" / 12pt run formatting used\n# throughout the document.\n\n$d = $word.ActiveDocument\n\n$newLines = @(\n    \"\",\n    \"Update: 2026-02-21 (Contacts)\",\n    \"- Added contact Segment mapping support (backend + frontend).\",\n    \"- Contacts now return segment name in list API for UI display.\",\n    \"- Added SegmentId to contact schema and upsert request.\",\n    \"- Outbound message send now auto-creates/updates contact by recipient number.\",\n    \"- Dashboard Contacts page now uses real segments API, supports create segment, and maps segment while creating contact.\",\n    \"- WABA Contacts page updated to map segment on add/edit and display segment column.\"\n)\n\nforeach ($line in $newLines) {\n    $tail = $d.Paragraphs.Last.Range\n    $tail.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Last.Range\n    $newPara.Text = $line\n    $newPara.Font.Name = \"Helvetica Light\"\n    $newPara.Font.Size = 12\n}\n"}
